# Updates sheet1 of the memory/living_rooms scene-category workbook:
# - block_total (C) 2 -> 4
# - trial_total (F) shifted/renumbered (81..108 -> 228..255) and 13 new trials appended (256..268)
# - category (H) kitchens -> bedrooms
# - cond_cat (I) recomputed per row (present only for "old"/target trials, absent for "new" trials)
# - cond_mem (J) / correct_answer (K) updated to reflect the new old/new split
# - stimulus (L) filenames updated
# - conceptual/perceptual/typicality/n/p_* (M:S) updated with the recomputed sanity-check numbers
# - 13 new trial rows (30:42) appended, increasing the proportion of "new" images in the memory task

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2, 'memory', 4, 1, 1, 228, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_6ddrx.png', 82.2, 63.68571428571428, 72.94285714285715, 35, 8, 8, 8),
    @(2, 'memory', 4, 1, 2, 229, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_b89t4.png', 71.41463414634147, 47.85365853658536, 59.63414634146342, 41, 5, 5, 5),
    @(2, 'memory', 4, 1, 3, 230, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_badai.png', 63.97435897435897, 43.38461538461539, 53.67948717948718, 39, 4, 4, 4),
    @(2, 'memory', 4, 1, 4, 231, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_0nckg.png', 65.94285714285714, 41.17142857142857, 53.55714285714285, 35, 4, 4, 4),
    @(2, 'memory', 4, 1, 5, 232, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_fbihy.png', 44.39024390243902, 26.90243902439024, 35.64634146341464, 41, 2, 2, 2),
    @(2, 'memory', 4, 1, 6, 233, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_ybbmx.png', 55.24324324324324, 36.75675675675676, 46, 37, 3, 3, 3),
    @(2, 'memory', 4, 1, 7, 234, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_0eflx.png', 76.05128205128206, 53.53846153846154, 64.7948717948718, 39, 6, 6, 6),
    @(2, 'memory', 4, 1, 8, 235, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_gqy6z.png', 86.47368421052632, 68.42105263157895, 77.44736842105263, 38, 9, 9, 9),
    @(2, 'memory', 4, 1, 9, 236, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_h0px5.png', 84.30232558139535, 69.13953488372093, 76.72093023255815, 43, 9, 9, 9),
    @(2, 'memory', 4, 1, 10, 237, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_qgbyn.png', 65.08108108108108, 40.10810810810811, 52.5945945945946, 37, 4, 4, 4),
    @(2, 'memory', 4, 1, 11, 238, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_wyctg.png', 33.44736842105263, 11.39473684210526, 22.42105263157895, 38, 1, 1, 1),
    @(2, 'memory', 4, 1, 12, 239, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_5p2ql.png', 89.19565217391305, 72.52173913043478, 80.8586956521739, 46, 10, 10, 10),
    @(2, 'memory', 4, 1, 13, 240, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_bj2gr.png', 65.25, 44.8, 55.025, 40, 4, 4, 4),
    @(2, 'memory', 4, 1, 14, 241, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_dmjh8.png', 57.48648648648648, 37.64864864864865, 47.56756756756756, 37, 3, 3, 3),
    @(2, 'memory', 4, 1, 15, 242, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_2dnfy.png', 59.97297297297298, 37.56756756756756, 48.77027027027027, 37, 3, 3, 3),
    @(2, 'memory', 4, 1, 16, 243, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_awuev.png', 44.21052631578947, 24.26315789473684, 34.23684210526316, 38, 2, 2, 2),
    @(2, 'memory', 4, 1, 17, 244, 'living_rooms', $null, $null, 'catch', 'f', 'stimuli/catch_03.jpg', $null, $null, $null, $null, $null, $null, $null),
    @(2, 'memory', 4, 1, 18, 245, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_71mhq.png', 69.34210526315789, 47.02631578947368, 58.18421052631579, 38, 5, 5, 5),
    @(2, 'memory', 4, 1, 19, 246, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_qbdgm.png', 76.88095238095238, 60.40476190476191, 68.64285714285714, 42, 7, 7, 7),
    @(2, 'memory', 4, 1, 20, 247, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_psgf7.png', 26, 11.66666666666667, 18.83333333333333, 36, 1, 1, 1),
    @(2, 'memory', 4, 1, 21, 248, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_iqmdm.png', 79.38888888888889, 58.36111111111111, 68.875, 36, 7, 7, 7),
    @(2, 'memory', 4, 1, 22, 249, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_q8u19.png', 91.1951219512195, 73.58536585365853, 82.39024390243901, 41, 10, 10, 10),
    @(2, 'memory', 4, 1, 23, 250, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_x0u5z.png', 92, 78.16216216216216, 85.08108108108108, 37, 10, 10, 10),
    @(2, 'memory', 4, 1, 24, 251, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_g2akb.png', 87.875, 79, 83.4375, 40, 10, 10, 10),
    @(2, 'memory', 4, 1, 25, 252, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_le8uf.png', 12.88888888888889, 9.222222222222221, 11.05555555555556, 36, 1, 1, 1),
    @(2, 'memory', 4, 1, 26, 253, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_uxxo0.png', 71.74418604651163, 48.44186046511628, 60.09302325581395, 43, 5, 5, 5),
    @(2, 'memory', 4, 1, 27, 254, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_oou46.png', 75.70270270270271, 54.86486486486486, 65.28378378378379, 37, 6, 6, 6),
    @(2, 'memory', 4, 1, 28, 255, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_5m6x4.png', 80.23076923076923, 58.41025641025641, 69.32051282051282, 39, 7, 7, 7),
    @(2, 'memory', 4, 1, 29, 256, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_h0hbk.png', 86.80952380952381, 69.19047619047619, 78, 42, 9, 9, 9),
    @(2, 'memory', 4, 1, 30, 257, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_okvvw.png', 50.58333333333334, 32.11111111111111, 41.34722222222223, 36, 2, 2, 2),
    @(2, 'memory', 4, 1, 31, 258, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_2js6m.png', 40.02777777777778, 20.88888888888889, 30.45833333333334, 36, 2, 2, 2),
    @(2, 'memory', 4, 1, 32, 259, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_rvssl.png', 74.25, 54.33333333333334, 64.29166666666667, 36, 6, 6, 6),
    @(2, 'memory', 4, 1, 33, 260, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_ozxpp.png', 26.26470588235294, 11.47058823529412, 18.86764705882353, 34, 1, 1, 1),
    @(2, 'memory', 4, 1, 34, 261, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_a8y4y.png', 75.15789473684211, 53.76315789473684, 64.46052631578948, 38, 6, 6, 6),
    @(2, 'memory', 4, 1, 35, 262, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_9z99v.png', 81.15625, 64.78125, 72.96875, 32, 8, 8, 8),
    @(2, 'memory', 4, 1, 36, 263, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_a9acb.png', 77.11428571428571, 58.42857142857143, 67.77142857142857, 35, 7, 7, 7),
    @(2, 'memory', 4, 1, 37, 264, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_v8dra.png', 61.77272727272727, 38.79545454545455, 50.28409090909091, 44, 3, 3, 3),
    @(2, 'memory', 4, 1, 38, 265, 'living_rooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_wijef.png', 69.875, 48.025, 58.95, 40, 5, 5, 5),
    @(2, 'memory', 4, 1, 39, 266, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_2pk6v.png', 85.08108108108108, 66.16216216216216, 75.62162162162161, 37, 9, 9, 9),
    @(2, 'memory', 4, 1, 40, 267, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_fqgem.png', 80.75, 61.475, 71.1125, 40, 8, 8, 8),
    @(2, 'memory', 4, 1, 41, 268, 'living_rooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_th7xh.png', 82.35897435897436, 65.53846153846153, 73.94871794871796, 39, 8, 8, 8)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    for ($c = 1; $c -le $vals.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}
